$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Production values (column C) for rows 2..170, in order.
$prod = @(
    0.0002668727829586715, 0.0003325073630549014, 0.2648947536945343, 0.3087421655654907, 0.3379668593406677, 0.3341384828090668, 0.3193267583847046, 0.31169193983078, 0.2951627969741821, 0.2422835677862167,
    0.1868152171373367, 0.1064927577972412, 0.0942598506808281, 0.03408302366733551, 0.001089186756871641, -0.00006791258056182414, -0.00006791258056182414, -0.0000903630570974201, -0.000009281180609832518, 0.0001480838982388377,
    0.0001480838982388377, 0.0001480838982388377, 0.0008348786504939198, 0.03718990832567215, 0.1236515566706657, 0.2028597891330719, 0.2627995908260345, 0.312704861164093, 0.3354771733283997, 0.3336235582828522,
    0.3190118968486786, 0.3151858747005463, 0.2991250455379486, 0.2424368560314178, 0.1862100064754486, 0.1034642606973648, 0.09533266723155975, 0.0332803875207901, -0.00001242202597495634, -0.0003874365356750786,
    0.0005688919336535037, 0.0005909791216254234, -0.00004582531619234942, -0.00004582531619234942, -0.00004582531619234942, -0.00004582531619234942, 0.0005531444330699742, 0.02721544727683067, 0.1087944805622101, 0.19118632376194,
    0.2583509981632233, 0.3022304475307465, 0.3244161903858185, 0.3282375335693359, 0.308084100484848, 0.2879491150379181, 0.2609406113624573, 0.2179675251245499, 0.150703638792038, 0.1016505211591721,
    0.08676796406507492, 0.0332803875207901, 0.0007015404989942908, -0.00006791258056182414, -0.00006791258056182414, -0.00004582531619234942, -0.0000903630570974201, -0.0000903630570974201, 0.00006700201629428193, 0.00006700201629428193,
    0.0007537968340329826, 0.03710882738232613, 0.1156279146671295, 0.1946479231119156, 0.2627995908260345, 0.3056575059890747, 0.3318667709827423, 0.3284536004066467, 0.3145216405391693, 0.3068116009235382,
    0.2832940816879272, 0.2355321943759918, 0.1738503277301788, 0.1017130464315414, 0.0910051241517067, 0.03100715391337872, -0.0008150548674166203, -0.0003874365356750786, -0.0003874365356750786, 0.001035457593388855,
    0.0005909791216254234, -0.00004582531619234942, -0.00004582531619234942, 0.00003525653664837591, 0.0006342263077385724, 0.03700832650065422, 0.1165561899542809, 0.19118632376194, 0.2602411508560181, 0.2999543249607086,
    0.3229856789112091, 0.3273255527019501, 0.316826730966568, 0.3074976205825806, 0.2772199213504791, 0.2161193042993546, 0.1491924524307251, 0.09208405762910843, 0.05837381258606911, 0.01870098896324635,
    -0.00003299403397249989, 0.0003088484518229961, -0.00005910069739911705, -0.00003701345121953636, -0.00003701345121953636, -0.00004582531619234942, -0.000009281180609832518, -0.000009281180609832518, 0.0006342263077385724, 0.03700832650065422,
    0.112518385052681, 0.19118632376194, 0.2594632804393768, 0.3033076822757721, 0.3244161903858185, 0.3212643265724182, 0.2945206165313721, 0.2532672584056854, 0.1942060589790344, 0.1381801962852478,
    0.1032142490148544, 0.08097250759601593, 0.07150858640670776, 0.02509656175971031, 0.0004092000308446586, -0.00008569909550715238, -0.0001036383328028023, -0.00008155114483088255, -0.00008155114483088255, -0.00008155114483088255,
    -0.00009946202044375241, -0.00009946202044375241, 0.0004945086548104882, 0.02715681120753288, 0.1164467036724091, 0.1908307671546936, 0.2594632804393768, 0.3022522032260895, 0.3282984793186188, 0.3282375335693359,
    0.314129650592804, 0.3059799373149872, 0.2836954891681671, 0.2359336167573929, 0.1779517531394958, 0.1014163345098495, 0.09966443479061127, 0.03206797316670418, -0.00001242202597495634, 0.0003088484518229961,
    -0.00005910069739911705, -0.00005492430864251219, -0.0001044611271936446, -0.0000903630570974201, -0.0000903630570974201, -0.0000903630570974201, 0.0005531444330699742, 0.03371951729059219, 0.112518385052681
)

# Date (column A, Excel serial) / Interval-hour (column B) blocks: (startRow, endRow, startDate, startHour)
$blocks = @(
    @(2, 18, 45461, 7),
    @(19, 42, 45462, 0),
    @(43, 66, 45463, 0),
    @(67, 90, 45464, 0),
    @(91, 114, 45465, 0),
    @(115, 138, 45466, 0),
    @(139, 162, 45467, 0),
    @(163, 170, 45468, 0)
)

foreach ($b in $blocks) {
    $startRow = $b[0]
    $endRow = $b[1]
    $date = $b[2]
    $startHour = $b[3]
    for ($r = $startRow; $r -le $endRow; $r++) {
        $hour = $startHour + ($r - $startRow)
        $ws.Cells.Item($r, 1).Value = $date
        $ws.Cells.Item($r, 2).Value = $hour
        $ws.Cells.Item($r, 3).Value = $prod[$r - 2]
    }
}

